$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.017.07"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.55%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.709.21"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.07%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9974"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.62%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "317.71"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.22%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.22%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4041"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +1.71%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4085"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.92%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.484"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -2.56%  "
$ws.Range("B10").Value = "OKB"
$ws.Range("C10").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "53.74"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.28%  "
$ws.Range("B11").Value = "BinanceUSD"
$ws.Range("C11").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9990"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.40%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08848"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.74%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "26.62"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +9.03%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.505"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -2.56%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.162"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.64%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001364"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.54%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.708.74"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.66%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "97.25"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -3.27%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.07169"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "21.27"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +5.08%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.285"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -2.92%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.000"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.75%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.41"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -1.24%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "25.000.05"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.49%  "
$ws.Range("B25").Value = "LidoDAOToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.938"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -5.19%  "
$ws.Range("B26").Value = "Toncoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.329"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.94%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "23.36"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.15%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.288"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +19.89%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "167.07"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.11%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "145.83"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +4.34%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.461"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -9.22%  "
$ws.Range("B32").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C32").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.897.92"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.51%  "
$ws.Range("B33").Value = "WEMIXTOKEN"
$ws.Range("C33").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.232"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +13.52%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08855"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -3.01%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.03206"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +5.34%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "7.254"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -7.98%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.036"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -5.33%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2880"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +1.63%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.8506"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +1.99%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "10.94"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -2.31%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.09344"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.22%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "14.22"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -3.15%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.469"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -1.50%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "17.44"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +4.67%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.724"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +2.32%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.7459"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.69%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.245"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.82%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.408"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +3.88%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.9994"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.34%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "141.87"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.34%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.08367"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +3.47%  "
